# ------------------------------------------------------------------
# Adds a new "CCGL" cohort to the results workbook:
#   - comparativo_master: 20 new per-person rows (rows 260-279)
#   - niveis_master:      3 new per-level rows for CCGL + updated TOTAL rows
#                          + an AutoFilter on Nível = "Intermediário"
#   - financeiro_master:  1 new CCGL total row + updated TOTAL row
#   - financeiro_master becomes the active sheet/tab
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsComparativo = $wb.Worksheets.Item("comparativo_master")
$wsNiveis      = $wb.Worksheets.Item("niveis_master")
$wsFinanceiro  = $wb.Worksheets.Item("financeiro_master")

# ------------------------------------------------------------------
# 1) comparativo_master — append the 20 CCGL rows (260-279)
# ------------------------------------------------------------------
$comparativoData = @(
    ,@('CCGL','ARNO KOCHHANN HECH',32,'Intermediário',88,'Avançado',56,175)
    ,@('CCGL','CARLOS ALBERTO ARNT',23,'Básico',49,'Intermediário',26,113.04347826086961)
    ,@('CCGL','DALMIR ANDERSON KEMPF',10,'Básico',88,'Avançado',78,780)
    ,@('CCGL','DIOVANI VILLANI LOPES',30,'Intermediário',30,'Intermediário',0,0)
    ,@('CCGL','DIRCE MARLI TOSSIN',42,'Intermediário',65,'Avançado',23,54.761904761904773)
    ,@('CCGL','ELEANDRO BERNARDI',11,'Básico',78,'Avançado',67,609.09090909090912)
    ,@('CCGL','ELENI REGINA GONZATTO',37,'Intermediário',67,'Avançado',30,81.081081081081081)
    ,@('CCGL','FERNANDO BUENO SIMOES PIRES',12,'Básico',84,'Avançado',72,600)
    ,@('CCGL','GILMAR BALBINOT',43,'Intermediário',83,'Avançado',40,93.023255813953483)
    ,@('CCGL','IRINEU ANTONIO ROHR',11,'Básico',81,'Avançado',70,636.36363636363637)
    ,@('CCGL','JOAO RAFAEL BONNE',16,'Básico',58,'Intermediário',42,262.5)
    ,@('CCGL','MARIA HELENA PASA SCHAEFER',13,'Básico',54,'Intermediário',41,315.38461538461542)
    ,@('CCGL','MARINÊS PEDROZO',10,'Básico',76,'Avançado',66,660)
    ,@('CCGL','MOACIR ZANIN',3,'Básico',60,'Avançado',57,1900)
    ,@('CCGL','NADI ANDREIA KEMPF',40,'Intermediário',86,'Avançado',46,115)
    ,@('CCGL','Odacir Nicolodi',37,'Intermediário',75,'Avançado',38,102.70270270270269)
    ,@('CCGL','ROQUE ALBERTO ANHOLETTO',15,'Básico',86,'Avançado',71,473.33333333333331)
    ,@('CCGL','RUDINEI BONORA',49,'Intermediário',86,'Avançado',37,75.510204081632651)
    ,@('CCGL','SIRLENE CARVALHO DOS SANTOS',20,'Básico',52,'Intermediário',32,160)
    ,@('CCGL','VERA LUCIA MADALOZZO WEBER',35,'Intermediário',70,'Avançado',35,100)
)

$startRow = 260
for ($i = 0; $i -lt $comparativoData.Count; $i++) {
    $row = $startRow + $i
    $d = $comparativoData[$i]
    $wsComparativo.Cells.Item($row,1).Value = $d[0]
    $wsComparativo.Cells.Item($row,2).Value = $d[1]
    $wsComparativo.Cells.Item($row,3).Value = $d[2]
    $wsComparativo.Cells.Item($row,4).Value = $d[3]
    $wsComparativo.Cells.Item($row,5).Value = $d[4]
    $wsComparativo.Cells.Item($row,6).Value = $d[5]
    $wsComparativo.Cells.Item($row,7).Value = $d[6]
    $wsComparativo.Cells.Item($row,8).Value = $d[7]
}

$wsComparativo.Range("A261:A279").Select()

# ------------------------------------------------------------------
# 2) niveis_master — insert the CCGL level rows ahead of the TOTAL
#    rows, refresh the TOTAL rows, and filter on "Intermediário"
# ------------------------------------------------------------------

# Drop any existing AutoFilter first so re-applying it below recomputes
# the filtered range (and the hidden rows) against the full new extent
# instead of re-using the stale A1:D43 range.
if ($wsNiveis.AutoFilterMode()) {
    $wsNiveis.AutoFilterMode = $false
}

# Insert 3 blank rows right before the old row 41 (TOTAL/Avançado),
# pushing the three TOTAL rows down to 44-46.
$wsNiveis.Rows.Item(41).Resize(3,1).Insert()

$wsNiveis.Cells.Item(41,1).Value = "CCGL"
$wsNiveis.Cells.Item(41,2).Value = "Básico"
$wsNiveis.Cells.Item(41,3).Value = 11
$wsNiveis.Cells.Item(41,4).Value = 0

$wsNiveis.Cells.Item(42,1).Value = "CCGL"
$wsNiveis.Cells.Item(42,2).Value = "Intermediário"
$wsNiveis.Cells.Item(42,3).Value = 9
$wsNiveis.Cells.Item(42,4).Value = 5

$wsNiveis.Cells.Item(43,1).Value = "CCGL"
$wsNiveis.Cells.Item(43,2).Value = "Avançado"
$wsNiveis.Cells.Item(43,3).Value = 0
$wsNiveis.Cells.Item(43,4).Value = 15

# Refreshed TOTAL rows (old values + the new CCGL contribution)
$wsNiveis.Cells.Item(44,3).Value = 2
$wsNiveis.Cells.Item(44,4).Value = 158

$wsNiveis.Cells.Item(45,3).Value = 218
$wsNiveis.Cells.Item(45,4).Value = 21

$wsNiveis.Cells.Item(46,3).Value = 58
$wsNiveis.Cells.Item(46,4).Value = 99

# Re-apply the AutoFilter over the full A1:D46 range, restricted to
# Nível ("Intermediário") — this both hides the non-matching rows and
# writes the discrete <filters><filter> list.
$wsNiveis.Range("A1:D46").AutoFilter(2, @("Intermediário"), 7)

# The hidden "_FilterDatabase" defined name needs to track the new range.
foreach ($n in $wb.Names) {
    if ($n.Name() -eq "niveis_master!_FilterDatabase") {
        $n.RefersTo = "=niveis_master!`$A`$1:`$D`$46"
    }
}

$wsNiveis.Range("D47").Select()

# ------------------------------------------------------------------
# 3) financeiro_master — insert the CCGL total row ahead of the TOTAL
#    row and refresh the TOTAL row
# ------------------------------------------------------------------
$wsFinanceiro.Rows.Item(15).Insert()

$wsFinanceiro.Cells.Item(15,1).Value = "CCGL"
$wsFinanceiro.Cells.Item(15,2).Value = "Gestão Financeira"
$wsFinanceiro.Cells.Item(15,3).Value = 36
$wsFinanceiro.Cells.Item(15,4).Value = 197
$wsFinanceiro.Cells.Item(15,5).Value = 161
$wsFinanceiro.Cells.Item(15,6).Value = 447.22222222222217

$wsFinanceiro.Cells.Item(16,3).Value = 446
$wsFinanceiro.Cells.Item(16,4).Value = 2759
$wsFinanceiro.Cells.Item(16,5).Value = 2297
$wsFinanceiro.Cells.Item(16,6).Value = 516.92913385826773

$wsFinanceiro.Range("F16").Select()

# financeiro_master becomes the active sheet/tab (was status_consultorias)
$wsFinanceiro.Activate()
